# Auto-generated script applying betting odds updates for 2026-01-19 workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 1.46
$ws.Range("G2").Value = 1.55
$ws.Range("I2").Value = 8.6
$ws.Range("W2").Value = 2.8
$ws.Range("AJ2").Value = 16.5
$ws.Range("AN2").Value = 7.6
# Row 3
$ws.Range("R3").Value = 1.48
# Row 4
$ws.Range("F4").Value = 1.59
$ws.Range("H4").Value = 6
$ws.Range("K4").Value = 4.5
$ws.Range("P4").Value = 2.04
$ws.Range("Q4").Value = 1.87
$ws.Range("R4").Value = 1.4
$ws.Range("S4").Value = 3.2
$ws.Range("T4").Value = 1.92
$ws.Range("U4").Value = 1.96
$ws.Range("Z4").Value = 65
$ws.Range("AA4").Value = 210
$ws.Range("AE4").Value = 110
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 100
$ws.Range("AL4").Value = 44
$ws.Range("AM4").Value = 140
$ws.Range("AO4").Value = 130
# Row 5
$ws.Range("F5").Value = 4.6
$ws.Range("G5").Value = 4.9
$ws.Range("H5").Value = 1.9
$ws.Range("I5").Value = 1.96
$ws.Range("J5").Value = 3.6
$ws.Range("K5").Value = 3.85
$ws.Range("P5").Value = 1.74
$ws.Range("Q5").Value = 2.16
$ws.Range("R5").Value = 1.28
$ws.Range("T5").Value = 2
$ws.Range("U5").Value = 1.85
$ws.Range("V5").Value = 2.04
$ws.Range("AL5").Value = 95
$ws.Range("AM5").Value = 170
$ws.Range("AN5").Value = 110
$ws.Range("AO5").Value = 18
# Row 6
$ws.Range("F6").Value = 1.41
$ws.Range("G6").Value = 1.52
$ws.Range("H6").Value = 9
$ws.Range("I6").Value = 12.5
$ws.Range("J6").Value = 4.2
$ws.Range("K6").Value = 5.5
$ws.Range("L6").Value = 1.42
$ws.Range("N6").Value = 3.2
$ws.Range("O6").Value = 1.36
$ws.Range("P6").Value = 1.75
$ws.Range("Q6").Value = 2
$ws.Range("R6").Value = 1.28
$ws.Range("S6").Value = 3.9
$ws.Range("T6").Value = 2.22
$ws.Range("U6").Value = 1.54
$ws.Range("W6").Value = 2.92
$ws.Range("Y6").Value = 34
$ws.Range("AB6").Value = 7.4
$ws.Range("AC6").Value = 13.5
$ws.Range("AE6").Value = 280
$ws.Range("AF6").Value = 8.6
$ws.Range("AJ6").Value = 14
$ws.Range("AN6").Value = 11
# Row 7
$ws.Range("N7").Value = 3.2
$ws.Range("O7").Value = 1.39
$ws.Range("Q7").Value = 2.14
$ws.Range("T7").Value = 1.93
$ws.Range("U7").Value = 1.93
$ws.Range("Z7").Value = 980
$ws.Range("AC7").Value = 8
# Row 8
$ws.Range("G8").Value = 1.43
$ws.Range("H8").Value = 7.8
$ws.Range("I8").Value = 8.199999999999999
$ws.Range("Q8").Value = 1.48
$ws.Range("R8").Value = 1.75
$ws.Range("S8").Value = 2.18
$ws.Range("T8").Value = 1.71
$ws.Range("U8").Value = 2.22
$ws.Range("V8").Value = 1.14
$ws.Range("AO8").Value = 85
# Row 9
$ws.Range("J9").Value = 3.15
$ws.Range("N9").Value = 3.85
$ws.Range("P9").Value = 1.99
$ws.Range("Q9").Value = 1.82
$ws.Range("S9").Value = 2.78
# Row 11
$ws.Range("I11").Value = 1.97
$ws.Range("R11").Value = 1.26
# Row 12
$ws.Range("I12").Value = 5
$ws.Range("J12").Value = 3.25
$ws.Range("N12").Value = 3.1
$ws.Range("O12").Value = 1.38
$ws.Range("Q12").Value = 2.12
$ws.Range("R12").Value = 1.27
$ws.Range("S12").Value = 3.9
$ws.Range("V12").Value = 1.25
$ws.Range("W12").Value = 1.76
$ws.Range("AF12").Value = 14.5
$ws.Range("AH12").Value = 980
# Row 13
$ws.Range("G13").Value = 2.8
$ws.Range("H13").Value = 3.25
$ws.Range("J13").Value = 2.86
$ws.Range("K13").Value = 3.8
$ws.Range("L13").Value = 1.48
$ws.Range("N13").Value = 2.94
$ws.Range("P13").Value = 1.66
$ws.Range("Q13").Value = 2.22
$ws.Range("R13").Value = 1.24
$ws.Range("S13").Value = 4.2
$ws.Range("T13").Value = 1.79
$ws.Range("U13").Value = 1.91
$ws.Range("V13").Value = 1.32
$ws.Range("X13").Value = 11.5
# Row 14
$ws.Range("F14").Value = 3.1
$ws.Range("G14").Value = 5.4
$ws.Range("I14").Value = 2.16
$ws.Range("J14").Value = 3.4
$ws.Range("R14").Value = 1.18
$ws.Range("S14").Value = 1.61
$ws.Range("V14").Value = 1.86
$ws.Range("W14").Value = 1.22
# Row 15
$ws.Range("F15").Value = 2.18
$ws.Range("G15").Value = 2.36
$ws.Range("H15").Value = 3.8
$ws.Range("M15").Value = 1.09
$ws.Range("Q15").Value = 2.28
$ws.Range("V15").Value = 1.31
# Row 16
$ws.Range("H16").Value = 3.35
$ws.Range("K16").Value = 5.8
$ws.Range("T16").Value = 1.42
$ws.Range("AB16").Value = 21
$ws.Range("AI16").Value = 32
# Row 17
$ws.Range("G17").Value = 4.7
$ws.Range("K17").Value = 5
$ws.Range("L17").Value = 1.26
# Row 18
$ws.Range("I18").Value = 6.2
$ws.Range("J18").Value = 3.05
$ws.Range("L18").Value = 1.42
$ws.Range("O18").Value = 1.39
$ws.Range("P18").Value = 1.68
$ws.Range("V18").Value = 1.19
# Row 19
$ws.Range("G19").Value = 2.64
$ws.Range("I19").Value = 3.55
$ws.Range("P19").Value = 1.57
$ws.Range("T19").Value = 2.04
$ws.Range("X19").Value = 9
$ws.Range("AJ19").Value = 40
# Row 20
$ws.Range("F20").Value = 2.66
$ws.Range("G20").Value = 2.92
$ws.Range("H20").Value = 2.78
$ws.Range("M20").Value = 1.09
$ws.Range("N20").Value = 3.1
$ws.Range("P20").Value = 1.71
$ws.Range("Q20").Value = 2.14
$ws.Range("R20").Value = 1.27
$ws.Range("U20").Value = 1.99
$ws.Range("AN20").Value = 40
# Row 21
$ws.Range("J21").Value = 3.05
$ws.Range("K21").Value = 3.1
$ws.Range("L21").Value = 1.54
$ws.Range("N21").Value = 2.96
$ws.Range("O21").Value = 1.49
$ws.Range("Q21").Value = 2.5
# Row 22
$ws.Range("F22").Value = 1.93
$ws.Range("G22").Value = 1.94
$ws.Range("J22").Value = 4.2
$ws.Range("K22").Value = 4.3
$ws.Range("P22").Value = 2.82
$ws.Range("Q22").Value = 1.54
$ws.Range("W22").Value = 2.06
$ws.Range("AJ22").Value = 22
# Row 23
$ws.Range("N23").Value = 3.35
$ws.Range("P23").Value = 1.82
$ws.Range("Q23").Value = 2.2
$ws.Range("R23").Value = 1.3
$ws.Range("S23").Value = 4.1
$ws.Range("U23").Value = 2.02
$ws.Range("X23").Value = 12
$ws.Range("AB23").Value = 9
$ws.Range("AM23").Value = 120
# Row 25
$ws.Range("F25").Value = 2.28
$ws.Range("H25").Value = 3.55
$ws.Range("I25").Value = 4
$ws.Range("J25").Value = 3.1
$ws.Range("K25").Value = 3.45
$ws.Range("L25").Value = 1.41
$ws.Range("N25").Value = 2.96
$ws.Range("O25").Value = 1.43
$ws.Range("P25").Value = 1.67
$ws.Range("Q25").Value = 2.18
$ws.Range("R25").Value = 1.25
$ws.Range("S25").Value = 3.8
$ws.Range("Z25").Value = 32
$ws.Range("AC25").Value = 9
$ws.Range("AK25").Value = 36
$ws.Range("AL25").Value = 60
$ws.Range("AN25").Value = 34
